$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row heights (points) for rows 1-32, matching the recalculated
# auto-fit heights after the font/content update.
$rowHeights = @{
    1  = 18.75
    2  = 18.75
    3  = 18.75
    4  = 18.75
    5  = 32.25
    6  = 18.75
    7  = 32.25
    8  = 32.25
    9  = 45
    10 = 18.75
    11 = 18.75
    12 = 18.75
    13 = 18.75
    14 = 18.75
    15 = 18.75
    16 = 32.25
    17 = 32.25
    18 = 45.75
    19 = 18.75
    20 = 18.75
    21 = 18.75
    22 = 18.75
    23 = 18.75
    24 = 18.75
    25 = 18.75
    26 = 18.75
    27 = 18.75
    28 = 18.75
    29 = 18.75
    30 = 18.75
    31 = 18.75
    32 = 18.75
}

foreach ($rowNum in $rowHeights.Keys) {
    $ws.Rows.Item($rowNum).RowHeight = $rowHeights[$rowNum]
}
